$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite rows 23-49 in place with the reordered / renamed content
# (A_EU and A_GERMANY moved from the end of the sheet to the front of
# the K_AREA/K_LAENDER block; everything else in that block keeps its
# existing formatting because the cells are simply overwritten, not
# inserted/deleted).

$ws.Range("A23").Value = "A_EU"
$ws.Range("B23").Value = "K_AREA"
$ws.Range("C23").Value = "Eropäische Union"
$ws.Range("D23").Value = "European Union"

$ws.Range("A24").Value = "A_GERMANY"
$ws.Range("B24").Value = "K_AREA"
$ws.Range("C24").Value = "Deutschland"
$ws.Range("D24").Value = "Germany"

$ws.Range("A25").Value = "A_LAENDER_BB"
$ws.Range("B25").Value = "K_LAENDER"
$ws.Range("C25").Value = "Brandenburg"
$ws.Range("D25").Value = "Brandenburg"

$ws.Range("A26").Value = "A_LAENDER_BE"
$ws.Range("B26").Value = "K_LAENDER"
$ws.Range("C26").Value = "Berlin"
$ws.Range("D26").Value = "Berlin"

$ws.Range("A27").Value = "A_LAENDER_BW"
$ws.Range("B27").Value = "K_LAENDER"
$ws.Range("C27").Value = "Baden-Württemberg"
$ws.Range("D27").Value = "Baden-Wuerttemberg"

$ws.Range("A28").Value = "A_LAENDER_BY"
$ws.Range("B28").Value = "K_LAENDER"
$ws.Range("C28").Value = "Bayern"
$ws.Range("D28").Value = "Bavaria"

$ws.Range("A29").Value = "A_LAENDER_HB"
$ws.Range("B29").Value = "K_LAENDER"
$ws.Range("C29").Value = "Bremen"
$ws.Range("D29").Value = "Bremen"

$ws.Range("A30").Value = "A_LAENDER_HE"
$ws.Range("B30").Value = "K_LAENDER"
$ws.Range("C30").Value = "Hessen"
$ws.Range("D30").Value = "Hesse"

$ws.Range("A31").Value = "A_LAENDER_HH"
$ws.Range("B31").Value = "K_LAENDER"
$ws.Range("C31").Value = "Hamburg"
$ws.Range("D31").Value = "Hamburg"

$ws.Range("A32").Value = "A_LAENDER_MV"
$ws.Range("B32").Value = "K_LAENDER"
$ws.Range("C32").Value = "Mecklenburg-Vorpommern"
$ws.Range("D32").Value = "Mecklenburg Western Pomerania"

$ws.Range("A33").Value = "A_LAENDER_NI"
$ws.Range("B33").Value = "K_LAENDER"
$ws.Range("C33").Value = "Niedersachsen"
$ws.Range("D33").Value = "Lower Saxony"

$ws.Range("A34").Value = "A_LAENDER_NW"
$ws.Range("B34").Value = "K_LAENDER"
$ws.Range("C34").Value = "Nordrhein-Westfalen"
$ws.Range("D34").Value = "North Rhine-Westphalia"

$ws.Range("A35").Value = "A_LAENDER_RP"
$ws.Range("B35").Value = "K_LAENDER"
$ws.Range("C35").Value = "Rheinland-Pfalz"
$ws.Range("D35").Value = "Rhineland Palatinate"

$ws.Range("A36").Value = "A_LAENDER_SH"
$ws.Range("B36").Value = "K_LAENDER"
$ws.Range("C36").Value = "Schleswig-Holstein"
$ws.Range("D36").Value = "Schleswig-Holstein"

$ws.Range("A37").Value = "A_LAENDER_SL"
$ws.Range("B37").Value = "K_LAENDER"
$ws.Range("C37").Value = "Saarland"
$ws.Range("D37").Value = "Saarland"

$ws.Range("A38").Value = "A_LAENDER_SN"
$ws.Range("B38").Value = "K_LAENDER"
$ws.Range("C38").Value = "Sachsen"
$ws.Range("D38").Value = "Saxony"

$ws.Range("A39").Value = "A_LAENDER_ST"
$ws.Range("B39").Value = "K_LAENDER"
$ws.Range("C39").Value = "Sachsen-Anhalt"
$ws.Range("D39").Value = "Saxony-Anhalt"

$ws.Range("A40").Value = "A_LAENDER_TH"
$ws.Range("B40").Value = "K_LAENDER"
$ws.Range("C40").Value = "Thüringen"
$ws.Range("D40").Value = "Thuringia"

$ws.Range("A41").Value = "A_PM2.5"
$ws.Range("B41").Value = "K_PM"
$ws.Range("C41").Value = "PM2.5"
$ws.Range("D41").Value = "PM2.5"

$ws.Range("A42").Value = "A_SEA_B"
$ws.Range("B42").Value = "K_SEA"
$ws.Range("C42").Value = "Ostsee"
$ws.Range("D42").Value = "Baltic sea"

$ws.Range("A43").Value = "A_SEA_N"
$ws.Range("B43").Value = "K_SEA"
$ws.Range("C43").Value = "Nordsee"
$ws.Range("D43").Value = "Greater North Sea"

$ws.Range("A44").Value = "A_SERIES_5YAVERAGE"
$ws.Range("B44").Value = "K_SERIES"
$ws.Range("C44").Value = "Gleitender Fünfjahresdurchschnitt mit Bezug auf das mittlere Jahr"
$ws.Range("D44").Value = "Moving five-year average shown for each middle year"

$ws.Range("A45").Value = "A_SERIES_ANNUALVAL"
$ws.Range("B45").Value = "K_SERIES"
$ws.Range("C45").Value = "Berechnete jährliche Werte"
$ws.Range("D45").Value = "Calculated annual values"

$ws.Range("A46").Value = "A_SERIES_BMEL"
$ws.Range("B46").Value = "K_SERIES"
$ws.Range("C46").Value = "Daten des BMEL"
$ws.Range("D46").Value = "Data from the Federal Ministry of Food and Agriculture"

$ws.Range("A47").Value = "A_SERIES_DSTTS"
$ws.Range("B47").Value = "K_SERIES"
$ws.Range("C47").Value = "Daten des Statistischen Bundesamtes"
$ws.Range("D47").Value = "Data from the Federal Statistical Office"

$ws.Range("A48").Value = "A_SEX_D"
$ws.Range("B48").Value = "K_SEX"
$ws.Range("C48").Value = "Divers"
$ws.Range("D48").Value = "Divers"

$ws.Range("A49").Value = "A_SEX_F"
$ws.Range("B49").Value = "K_SEX"
$ws.Range("C49").Value = "Weiblich"
$ws.Range("D49").Value = "Female"

# Rows 50-53 are brand new; give them the same formatting as the last
# existing data row before writing their values.
$ws.Range("A49:D49").Copy()
$ws.Range("A50:D53").PasteSpecial(-4122)

$ws.Range("A50").Value = "A_SEX_M"
$ws.Range("B50").Value = "K_SEX"
$ws.Range("C50").Value = "Männlich"
$ws.Range("D50").Value = "Male"

$ws.Range("A51").Value = "A_SEX_U"
$ws.Range("B51").Value = "K_SEX"
$ws.Range("C51").Value = "Unbekannt"
$ws.Range("D51").Value = "Unknown"

$ws.Range("A52").Value = "A_URBAN_NONRURAL"
$ws.Range("B52").Value = "K_URBAN"
$ws.Range("C52").Value = "Nicht-ländliche Gebiete"
$ws.Range("D52").Value = "Non-rural areas"

$ws.Range("A53").Value = "A_URBAN_RURAL"
$ws.Range("B53").Value = "K_URBAN"
$ws.Range("C53").Value = "Ländliche Gebiete"
$ws.Range("D53").Value = "Rural areas"
